$wb = $excel.ActiveWorkbook

# --- "runs" sheet updates ---
$runs = $wb.Worksheets.Item("runs")
$runs.Range("B1").Value = 1
$runs.Range("B2").Value = 15
$runs.Range("B3").Value = 100
$runs.Range("B3").Select()

# --- "params" sheet updates ---
$params = $wb.Worksheets.Item("params")
$params.Range("D6").Value = 0
$params.Range("D13").Select()

# --- Make "runs" the active/selected tab ---
$runs.Activate()
$runs.Select()
